$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Update the F-column (time_taken) timestamps on the "data" sheet ---
$ws1.Cells.Item(2, 6).Value = "2021-10-05 14:22:48.864846"
$ws1.Cells.Item(3, 6).Value = "2021-10-05 14:22:48.864855"
$ws1.Cells.Item(4, 6).Value = "2021-10-05 14:22:48.864859"
$ws1.Cells.Item(5, 6).Value = "2021-10-05 14:22:48.864861"
$ws1.Cells.Item(6, 6).Value = "2021-10-05 14:22:48.864864"
$ws1.Cells.Item(7, 6).Value = "2021-10-05 14:22:48.864867"
$ws1.Cells.Item(8, 6).Value = "2021-10-05 14:22:48.864870"
$ws1.Cells.Item(9, 6).Value = "2021-10-05 14:22:48.864873"
$ws1.Cells.Item(10, 6).Value = "2021-10-05 14:22:48.864876"
$ws1.Cells.Item(11, 6).Value = "2021-10-05 14:22:48.864879"
$ws1.Cells.Item(12, 6).Value = "2021-10-05 14:22:48.864881"
$ws1.Cells.Item(13, 6).Value = "2021-10-05 14:22:48.864884"
$ws1.Cells.Item(14, 6).Value = "2021-10-05 14:22:48.864887"
$ws1.Cells.Item(15, 6).Value = "2021-10-05 14:22:48.864890"
$ws1.Cells.Item(16, 6).Value = "2021-10-05 14:22:48.864892"
$ws1.Cells.Item(17, 6).Value = "2021-10-05 14:22:48.864895"
$ws1.Cells.Item(18, 6).Value = "2021-10-05 14:22:48.864898"
$ws1.Cells.Item(19, 6).Value = "2021-10-05 14:22:48.864901"
$ws1.Cells.Item(20, 6).Value = "2021-10-05 14:22:48.864903"
$ws1.Cells.Item(21, 6).Value = "2021-10-05 14:22:48.864906"
$ws1.Cells.Item(22, 6).Value = "2021-10-05 14:22:48.864909"
$ws1.Cells.Item(23, 6).Value = "2021-10-05 14:22:48.864911"

# --- 2. Add the new "metadata" worksheet, after "data" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws2.Name = "metadata"

# Header row (B1:G1) - reuse the bold/bordered header style from the "data" sheet
$ws2.Cells.Item(1, 2).Value = "data_name"
$ws2.Cells.Item(1, 3).Value = "data_id"
$ws2.Cells.Item(1, 4).Value = "data_version"
$ws2.Cells.Item(1, 5).Value = "data_version_created"
$ws2.Cells.Item(1, 6).Value = "panel_query_time"
$ws2.Cells.Item(1, 7).Value = "panel_get_request"

$ws1.Range("B1:F1").Copy()
$ws2.Range("B1:F1").PasteSpecial(-4122)

$ws1.Range("B1").Copy()
$ws2.Range("G1").PasteSpecial(-4122)

# Data row 2
$ws2.Cells.Item(2, 1).Value = 0
$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

$ws2.Cells.Item(2, 2).Value = "Skeletal Muscle Channelopathies"
$ws2.Cells.Item(2, 3).Value = 229

# "1.37" must be stored as TEXT (not auto-coerced to a number): force the
# cell to text format while assigning, then drop back to the default style
# so no stray number-format style is left behind on the cell.
$dCell = $ws2.Cells.Item(2, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.37"
$dCell.Style = "Normal"

$ws2.Cells.Item(2, 5).Value = "2021-07-20T14:07:04.406448Z"
$ws2.Cells.Item(2, 6).Value = "2021-10-05 14:22:48.861333"
$ws2.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/229/?format=json"

$ws1.Select()
$ws1.Range("A1").Select()
